$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 143
$ws.Range("I4").Value = 143
$ws.Range("K4").Value = 143
$ws.Range("M4").Value = -29
$ws.Range("H5").Value = 104.73684
$ws.Range("I5").Value = 109.28571
$ws.Range("K5").Value = 109.28571
$ws.Range("M5").Value = 5.714290000000005
$ws.Range("H34").Value = 12065.5
$ws.Range("I34").Value = 1120.6666
$ws.Range("K34").Value = 1120.6666
$ws.Range("M34").Value = -917.6666
$ws.Range("H36").Value = 12065.5
$ws.Range("I36").Value = 1120.6666
$ws.Range("K36").Value = 1120.6666
$ws.Range("M36").Value = -405.6666
$ws.Range("H57").Value = 34500
$ws.Range("J57").Value = 34500
$ws.Range("L57").Value = 103500
$ws.Range("N57").Value = -104498
$ws.Range("H132").Value = 2072.6606
$ws.Range("I132").Value = 1233.8163
$ws.Range("J132").Value = 7944.5713
$ws.Range("K132").Value = 3701.4489
$ws.Range("L132").Value = 23833.7139
$ws.Range("M132").Value = -1171.4489
$ws.Range("N132").Value = -28893.7139

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 300
$ws.Range("I4").Value = 300
$ws.Range("J4").Value = 300
$ws.Range("K4").Value = 300
$ws.Range("L4").Value = 300
$ws.Range("M4").Value = -184
$ws.Range("N4").Value = -532
$ws.Range("H39").Value = 13746
$ws.Range("I39").Value = 7829
$ws.Range("K39").Value = 7829
$ws.Range("M39").Value = -7309
$ws.Range("H61").Value = 1106.4
$ws.Range("I61").Value = 1058.3448
$ws.Range("K61").Value = 1058.3448
$ws.Range("M61").Value = -846.3448000000001
$ws.Range("H74").Value = 1836.85
$ws.Range("I74").Value = 740.41174
$ws.Range("J74").Value = 8050
$ws.Range("K74").Value = 740.41174
$ws.Range("L74").Value = 8050
$ws.Range("M74").Value = 133.58826
$ws.Range("N74").Value = -9798
$ws.Range("H77").Value = 1836.85
$ws.Range("I77").Value = 740.41174
$ws.Range("J77").Value = 8050
$ws.Range("K77").Value = 3702.0587
$ws.Range("L77").Value = 40250
$ws.Range("M77").Value = 665.9413
$ws.Range("N77").Value = -48986
$ws.Range("H132").Value = 38465268
$ws.Range("I132").Value = 50001910
$ws.Range("J132").Value = 9790.333000000001
$ws.Range("K132").Value = 150005730
$ws.Range("L132").Value = 29370.999
$ws.Range("M132").Value = -150003200
$ws.Range("N132").Value = -34430.999
$ws.Range("H136").Value = 1106.4
$ws.Range("I136").Value = 1058.3448
$ws.Range("K136").Value = 3175.0344
$ws.Range("M136").Value = -625.0344000000005
$ws.Range("H138").Value = 53409.168
$ws.Range("J138").Value = 53409.168
$ws.Range("L138").Value = 53409.168
$ws.Range("N138").Value = -63689.168

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H74").Value = 43406.715
$ws.Range("J74").Value = 43406.715
$ws.Range("L74").Value = 43406.715
$ws.Range("N74").Value = -45278.715
$ws.Range("H77").Value = 43406.715
$ws.Range("J77").Value = 43406.715
$ws.Range("L77").Value = 130220.145
$ws.Range("N77").Value = -139580.145
$ws.Range("H134").Value = 3010.1904
$ws.Range("I134").Value = 3166.3076
$ws.Range("J134").Value = 2756.5
$ws.Range("K134").Value = 9498.9228
$ws.Range("L134").Value = 8269.5
$ws.Range("M134").Value = -6963.9228
$ws.Range("N134").Value = -13339.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 30.583334
$ws.Range("I7").Value = 28.818182
$ws.Range("J7").Value = 50
$ws.Range("K7").Value = 28.818182
$ws.Range("L7").Value = 50
$ws.Range("M7").Value = 84.18181799999999
$ws.Range("N7").Value = -276
$ws.Range("H31").Value = 3268.2144
$ws.Range("I31").Value = 1482.4348
$ws.Range("K31").Value = 1482.4348
$ws.Range("M31").Value = -1187.4348
$ws.Range("H34").Value = 3268.2144
$ws.Range("I34").Value = 1482.4348
$ws.Range("K34").Value = 1482.4348
$ws.Range("M34").Value = -1280.4348
$ws.Range("H58").Value = 1382.6774
$ws.Range("I58").Value = 1516.0555
$ws.Range("J58").Value = 1198
$ws.Range("K58").Value = 1516.0555
$ws.Range("L58").Value = 1198
$ws.Range("M58").Value = -1313.0555
$ws.Range("N58").Value = -1604
$ws.Range("H132").Value = 2308.6667
$ws.Range("I132").Value = 2638.3635
$ws.Range("J132").Value = 1402
$ws.Range("K132").Value = 7915.0905
$ws.Range("L132").Value = 4206
$ws.Range("M132").Value = -5385.0905
$ws.Range("N132").Value = -9266
$ws.Range("H134").Value = 6408.1934
$ws.Range("I134").Value = 7016.815
$ws.Range("K134").Value = 21050.445
$ws.Range("M134").Value = -18515.445
$ws.Range("H136").Value = 1382.6774
$ws.Range("I136").Value = 1516.0555
$ws.Range("J136").Value = 1198
$ws.Range("K136").Value = 4548.166499999999
$ws.Range("L136").Value = 3594
$ws.Range("M136").Value = -1998.166499999999
$ws.Range("N136").Value = -8694
$ws.Range("H139").Value = 54332
$ws.Range("J139").Value = 54332
$ws.Range("L139").Value = 54332
$ws.Range("N139").Value = -64612

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 2323
$ws.Range("I34").Value = 385
$ws.Range("J34").Value = 2599.8572
$ws.Range("K34").Value = 1155
$ws.Range("L34").Value = 7799.571599999999
$ws.Range("M34").Value = -1071
$ws.Range("N34").Value = -7967.571599999999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3254.2354
$ws.Range("I132").Value = 3042.4
$ws.Range("J132").Value = 3556.8572
$ws.Range("K132").Value = 9127.200000000001
$ws.Range("L132").Value = 10670.5716
$ws.Range("M132").Value = -6597.200000000001
$ws.Range("N132").Value = -15730.5716
$ws.Range("H140").Value = 71853.336
$ws.Range("J140").Value = 71853.336
$ws.Range("L140").Value = 71853.336
$ws.Range("N140").Value = -82213.336

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3587.1333
$ws.Range("I132").Value = 2710
$ws.Range("J132").Value = 5999.25
$ws.Range("K132").Value = 8130
$ws.Range("L132").Value = 17997.75
$ws.Range("M132").Value = -5600
$ws.Range("N132").Value = -23057.75
$ws.Range("H136").Value = 1621.4419
$ws.Range("I136").Value = 1621.8948
$ws.Range("J136").Value = 1618
$ws.Range("K136").Value = 4865.6844
$ws.Range("L136").Value = 4854
$ws.Range("M136").Value = -2315.6844
$ws.Range("N136").Value = -9954

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3391.739
$ws.Range("I132").Value = 6129.6665
$ws.Range("J132").Value = 1631.6428
$ws.Range("K132").Value = 18388.9995
$ws.Range("L132").Value = 4894.928400000001
$ws.Range("M132").Value = -15858.9995
$ws.Range("N132").Value = -9954.928400000001
$ws.Range("H136").Value = 773
$ws.Range("I136").Value = 776.4909
$ws.Range("K136").Value = 2329.4727
$ws.Range("M136").Value = 220.5272999999997
